# Update the date line at the top of the document.
$d = $word.ActiveDocument
$d.Content.Find.Execute("2024-07-23 Tuesday", $true, $false, $false, $false, $false,
                         $true, 1, $false, "2024-07-24 Wednesday", 2)

# Update the multiplication-fact table cells.
# The table has 20 rows x 5 columns; only rows 1, 5, 10, 15, 20 contain text.
$t = $d.Tables.Item(1)

function Set-CellText($table, $row, $col, $newText) {
    $cell = $table.Cell($row, $col)
    $rng = $cell.Range
    # Exclude the trailing end-of-cell marker (2 chars: cr + bell) from the range.
    $rng.SetRange($rng.Start, $rng.End - 1)
    $rng.Text = $newText
}

Set-CellText $t 1 1 "159×7=1113"
Set-CellText $t 1 2 "967×9=8703"
Set-CellText $t 1 3 "641×4=2564"
Set-CellText $t 1 4 "816×9=7344"
Set-CellText $t 1 5 "167×9=1503"

Set-CellText $t 5 1 "857×9=7713"
Set-CellText $t 5 2 "843×8=6744"
Set-CellText $t 5 3 "704×4=2816"
Set-CellText $t 5 4 "633×7=4431"
Set-CellText $t 5 5 "833×3=2499"

Set-CellText $t 10 1 "488×6=2928"
Set-CellText $t 10 2 "818×2=1636"
Set-CellText $t 10 3 "356×9=3204"
Set-CellText $t 10 4 "469×6=2814"
Set-CellText $t 10 5 "432×5=2160"

Set-CellText $t 15 1 "427×4=1708"
Set-CellText $t 15 2 "632×9=5688"
Set-CellText $t 15 3 "224×5=1120"
Set-CellText $t 15 4 "598×7=4186"
Set-CellText $t 15 5 "904×5=4520"

Set-CellText $t 20 1 "288×8=2304"
Set-CellText $t 20 2 "765×9=6885"
Set-CellText $t 20 3 "427×2=854"
Set-CellText $t 20 4 "780×5=3900"
Set-CellText $t 20 5 "904×5=4520"
